$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1. Remove the "TS_001" worksheet entirely (kept only "Test Scenarios")
# ------------------------------------------------------------------
$wb.Worksheets("TS_001").Delete()

$ws = $wb.Worksheets("Test Scenarios")

# ------------------------------------------------------------------
# 2. Add three new rows (5, 6, 7) re-using the same look & feel as
#    row 4 (copy formatting only, then fill in the values).
# ------------------------------------------------------------------
$ws.Range("B4:F4").Copy()
$ws.Range("B5:F5").PasteSpecial(-4122)
$ws.Range("B4:F4").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)
$ws.Range("B4:F4").Copy()
$ws.Range("B7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 5 - "(TS_002) Login functionality"
$ws.Range("B5").Value = "(TS_002)`nLogin functionality"
$ws.Range("D5").Value = "Validate the working of Login account functionality"
$ws.Range("C5").Value = "FRS"
$ws.Range("E5").Value = "P0"
$ws.Range("F5").Value = 7

# Row 6 - "(TS_003) Home Page Currency"  (note: description typed before the title)
$ws.Range("D6").Value = "Validate the working of home page > Currency"
$ws.Range("B6").Value = "(TS_003)`nHome Page Currency"
$ws.Range("C6").Value = "FRS"
$ws.Range("E6").Value = "P0"
$ws.Range("F6").Value = 4

# Row 7 - "(TS_004) Home Page Contact Us"
$ws.Range("B7").Value = "(TS_004)`nHome Page Contact Us"
$ws.Range("D7").Value = "Validate the working of home page > Contact Us"
$ws.Range("C7").Value = "FRS"
$ws.Range("E7").Value = "P0"
$ws.Range("F7").ClearContents()

# Match the row height used by row 4 for all the freshly added rows
$ws.Rows("5").RowHeight = 30
$ws.Rows("6").RowHeight = 30
$ws.Rows("7").RowHeight = 30

# ------------------------------------------------------------------
# 3. Widen column F a bit (it no longer needs to "best fit" the old
#    header text, it now needs to comfortably fit under the new data)
# ------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 25.42

# ------------------------------------------------------------------
# 4. Final selection, like the author left the cursor on E7
# ------------------------------------------------------------------
$ws.Range("E7").Select()

Write-Host "done"
